$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A (row numbers 8,11,12,13,14,15 with a header-style
# formatting carried over) is removed entirely, and the remaining columns
# B:F shift left to become A:E. Deleting the entire column A reproduces
# exactly this: all data in B:F moves into A:E, and the now-unused column F
# disappears.
$ws.Columns.Item(1).Delete()
